# Apply crypto price/volume updates per diff (GitHub Actions refresh, 2024-02-01).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.248.76'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '2.269.93'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '300.11'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").Value = '95.91'
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.494'
$ws.Range("E7").Value = '  -2.33%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -1.77%  '
$ws.Range("D10").Value = '33.18'
$ws.Range("E10").Value = '  -3.81%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '48.41'
$ws.Range("E12").Value = '  -6.52%  '
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = '6.65'
$ws.Range("E14").Value = '  -1.25%  '
$ws.Range("D15").Value = '15.65'
$ws.Range("E15").Value = '  -0.29%  '
$ws.Range("D16").Value = '2.621.35'
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").Value = '2.283.51'
$ws.Range("E17").Value = '  -2.57%  '
$ws.Range("D18").Value = '0.783'
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("D19").Value = '42.159.85'
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = '11.65'
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("E22").Value = '  -0.90%  '
$ws.Range("D23").Value = '66.24'
$ws.Range("E23").Value = '  -2.33%  '
$ws.Range("D24").Value = '235.07'
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.90'
$ws.Range("E28").Value = '  -4.24%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '167.81'
$ws.Range("E30").Value = '  +2.62%  '
$ws.Range("D31").Value = '9.18'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '33.66'
$ws.Range("E32").Value = '  -3.07%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '4.88'
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("D35").Value = '4.56'
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("D36").Value = '16.72'
$ws.Range("E36").Value = '  -1.17%  '
$ws.Range("E37").Value = '  -3.43%  '
$ws.Range("D38").Value = '0.0684'
$ws.Range("E38").Value = '  -3.74%  '
$ws.Range("D39").Value = '2.79'
$ws.Range("E39").Value = '  -3.04%  '
$ws.Range("D40").Value = '0.0984'
$ws.Range("E40").Value = '  -2.10%  '
$ws.Range("D41").Value = '0.109'
$ws.Range("E41").Value = '  -2.49%  '
$ws.Range("D42").Value = '1.72'
$ws.Range("E42").Value = '  -4.50%  '
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("D44").Value = '1.968.50'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("D45").Value = '0.0277'
$ws.Range("E45").Value = '  -1.02%  '
$ws.Range("D46").Value = '17.43'
$ws.Range("E46").Value = '  -5.96%  '
$ws.Range("D47").Value = '9.53'
$ws.Range("E47").Value = '  -6.36%  '
$ws.Range("E48").Value = '  -4.70%  '
$ws.Range("D49").Value = '2.493.71'
$ws.Range("E49").Value = '  -1.04%  '
$ws.Range("D50").Value = '52.28'
$ws.Range("E50").Value = '  -5.71%  '
$ws.Range("D51").Value = '1.48'
$ws.Range("E51").Value = '  -0.65%  '
